$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Update "numeroCuenta" column (Q) for rows 2 and 3: 406-740100-05 -> 406-785280-05 ---
# (Done first so the new shared string for the account number is appended before
#  the new shared string for the username, matching insertion order.)
$ws.Range("Q2").Value = "406-785280-05"
$ws.Range("Q3").Value = "406-785280-05"

# --- Update "usuario" column (D) for rows 2 and 3: userqa10 -> ospciclo4finde ---
# A plain Value assignment on these cells resets their cell style, so the
# existing formatting is preserved via a copy / paste-special(formats) round trip.
$tempRange = $ws.Range("Z100")

$ws.Range("D2").Copy($tempRange)
$ws.Range("D2").Value = "ospciclo4finde"
$tempRange.Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("D3").Copy($tempRange)
$ws.Range("D3").Value = "ospciclo4finde"
$tempRange.Copy()
$ws.Range("D3").PasteSpecial(-4122)

$tempRange.Clear()
$excel.CutCopyMode = 0

# --- Update sheet view: clear the frozen/top-left cell and move the selection to G7 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G7").Select()

$wb.Save()
